$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 4 new columns before the existing "molfile" column (old X, now will be AB)
# so headers shift: old X (molfile) -> AB, old Y (canonical smiles) -> AC
$ws.Range("X1:AA1").EntireColumn.Insert()

# Populate the new header cells with the new field names
$ws.Range("X1").Value = "color"
$ws.Range("Y1").Value = "solubility"
$ws.Range("Z1").Value = "form"
$ws.Range("AA1").Value = "inventory label"

# Match the header cell style used by the rest of row 1 (style index 1)
$ws.Range("X1:AA1").Style = $ws.Range("W1").Style

# Column widths per updated layout
$ws.Columns.Item(27).ColumnWidth = 18.28515625

# Update selection / view to match
$ws.Range("Z5").Select()
